$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 46
$ws1.Range("F5").Value = 520
$ws1.Range("F6").Value = 1561
$ws1.Range("F8").Value = 1208
$ws1.Range("F10").Value = 220
$ws1.Range("F11").Value = 162
$ws1.Range("F17").Value = 193
$ws1.Range("F18").Value = 180

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 5

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 46
$ws4.Range("F5").Value = 520
$ws4.Range("F6").Value = 1561
$ws4.Range("F7").Value = 5
$ws4.Range("F9").Value = 1209
$ws4.Range("F11").Value = 220
$ws4.Range("F12").Value = 162
$ws4.Range("F18").Value = 193
$ws4.Range("F19").Value = 180
